$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Locate the "MVPA min/week - Machine learning" row dynamically; the two new
# rows (MVPA_100 / MVPA_150) need to be inserted directly above it (i.e.
# directly below "MVPA min/week - ENMO, Median (Q1, Q3)").
$targetRow = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $c1 = $t.Cell($r, 1).Range.Text
    if ($c1 -like "*MVPA min/week - Machine learning*") {
        $targetRow = $r
        break
    }
}

if ($targetRow -eq 0) {
    Write-Output "ERROR: could not locate target row"
}
else {
    # Insert two new rows right before $targetRow, cloning formatting from it.
    $newRow1 = $t.Rows.Add($t.Rows.Item($targetRow))
    $newRow1.Cells.Item(1).Range.Text = "MVPA_100, Median (Q1, Q3)"
    $newRow1.Cells.Item(2).Range.Text = "707.8 (515.8, 936.6)"

    $newRow2 = $t.Rows.Add($t.Rows.Item($targetRow + 1))
    $newRow2.Cells.Item(1).Range.Text = "MVPA_150, Median (Q1, Q3)"
    $newRow2.Cells.Item(2).Range.Text = "304.8 (198.3, 443.5)"

    Write-Output ("Inserted rows before row " + $targetRow + "; table now has " + $t.Rows.Count + " rows")
}
